# The data rows (2-27) of the sheet are reshuffled: each row's content
# (columns A:T) is replaced by the content that another (original) row
# held before the edit. Row 1 (header) and row 2 keep their content.
#
# This mapping was derived by diffing the pre- and post-edit OOXML and
# matching each post-edit row to the pre-edit row whose values it now
# carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: snapshot the current ("before") values of every data row,
# columns A through T, before any writes happen.
$rows = @{}
for ($r = 2; $r -le 27; $r++) {
    $rows[$r] = $ws.Range("A" + $r + ":T" + $r).Value()
}

# Step 2: new row number -> source (original) row number.
$map = @{
    2  = 2
    3  = 26
    4  = 18
    5  = 12
    6  = 27
    7  = 19
    8  = 20
    9  = 8
    10 = 9
    11 = 14
    12 = 15
    13 = 16
    14 = 13
    15 = 3
    16 = 17
    17 = 5
    18 = 6
    19 = 4
    20 = 11
    21 = 22
    22 = 21
    23 = 23
    24 = 10
    25 = 7
    26 = 25
    27 = 24
}

# Step 3: write each row's new content from the captured snapshot.
foreach ($r in $map.Keys) {
    $src = $map[$r]
    $ws.Range("A" + $r + ":T" + $r).Value = $rows[$src]
}
